$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update status cells from "进行中"/"未开始" to "已完成"
$ws.Range("C33").Value = "已完成"
$ws.Range("C34").Value = "已完成"
$ws.Range("C35").Value = "已完成"
$ws.Range("C36").Value = "已完成"
$ws.Range("C37").Value = "已完成"
$ws.Range("C38").Value = "已完成"

# Update the summary cell text (merged A39:D40)
$ws.Range("A39").Value = "总结：此阶段做了以上界面的设计初稿，以图片形式插入的UI设计文档中。由于时间关系，设计稿附属文字描述暂无，后边会补上。"

# Update the selected cell in the sheet view
$ws.Range("C44").Select()
